$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme
for ($i=1; $i -le $cs.Count; $i++) {
  $c = $cs.Item($i)
  $rgb = $c.RGB
  $hex = "{0:X6}" -f $rgb
  Write-Output "$i : $rgb  hex=$hex"
}
